# Adds 4 new match rows (106-109) to Sheet1, mirroring the style/format
# of the existing data rows (row 105 for column A index style,
# and the date column E style), matching the source script update
# dated 26-11-2023 20:30.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 105   # last pre-existing data row, used as style template

# ---- Row 106 ----
$ws.Cells.Item(106,1).Value = 105
$ws.Cells.Item(106,2).Value = 'ecuador'
$ws.Cells.Item(106,3).Value = 'liga-pro'
$ws.Cells.Item(106,4).NumberFormat = "@"
$ws.Cells.Item(106,4).Value = '2023'
$ws.Cells.Item(106,5).Value = 45256
$ws.Cells.Item(106,6).Value = 'Aucas'
$ws.Cells.Item(106,7).Value = 0
$ws.Cells.Item(106,8).Value = 'Delfin'
$ws.Cells.Item(106,9).Value = 0
$ws.Cells.Item(106,10).Value = 1.71
$ws.Cells.Item(106,11).Value = '19/11/2023 00:12'
$ws.Cells.Item(106,12).Value = 1.89
$ws.Cells.Item(106,13).Value = '25/11/2023 23:59'
$ws.Cells.Item(106,14).Value = 3.64
$ws.Cells.Item(106,15).Value = '19/11/2023 00:12'
$ws.Cells.Item(106,16).Value = 3.59
$ws.Cells.Item(106,17).Value = '25/11/2023 23:59'
$ws.Cells.Item(106,18).Value = 4.59
$ws.Cells.Item(106,19).Value = '19/11/2023 00:12'
$ws.Cells.Item(106,20).Value = 4.17
$ws.Cells.Item(106,21).Value = '25/11/2023 23:52'
$ws.Cells.Item(106,22).Value = 'https://www.betexplorer.com/football/ecuador/liga-pro/aucas-delfin/vBaASdYr/'
$ws.Cells.Item($lastDataRow,1).Copy() | Out-Null
$ws.Cells.Item(106,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item($lastDataRow,5).Copy() | Out-Null
$ws.Cells.Item(106,5).PasteSpecial(-4122) | Out-Null

# ---- Row 107 ----
$ws.Cells.Item(107,1).Value = 106
$ws.Cells.Item(107,2).Value = 'ecuador'
$ws.Cells.Item(107,3).Value = 'liga-pro'
$ws.Cells.Item(107,4).NumberFormat = "@"
$ws.Cells.Item(107,4).Value = '2023'
$ws.Cells.Item(107,5).Value = 45256
$ws.Cells.Item(107,6).Value = 'Dep. Cuenca'
$ws.Cells.Item(107,7).Value = 1
$ws.Cells.Item(107,8).Value = 'EL Nacional'
$ws.Cells.Item(107,9).Value = 0
$ws.Cells.Item(107,10).Value = 2.33
$ws.Cells.Item(107,11).Value = '19/11/2023 00:12'
$ws.Cells.Item(107,12).Value = 3.04
$ws.Cells.Item(107,13).Value = '25/11/2023 23:55'
$ws.Cells.Item(107,14).Value = 3.44
$ws.Cells.Item(107,15).Value = '19/11/2023 00:12'
$ws.Cells.Item(107,16).Value = 3.39
$ws.Cells.Item(107,17).Value = '25/11/2023 23:55'
$ws.Cells.Item(107,18).Value = 2.98
$ws.Cells.Item(107,19).Value = '19/11/2023 00:12'
$ws.Cells.Item(107,20).Value = 2.38
$ws.Cells.Item(107,21).Value = '25/11/2023 23:55'
$ws.Cells.Item(107,22).Value = 'https://www.betexplorer.com/football/ecuador/liga-pro/dep-cuenca-el-nacional/x2DrXxQR/'
$ws.Cells.Item($lastDataRow,1).Copy() | Out-Null
$ws.Cells.Item(107,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item($lastDataRow,5).Copy() | Out-Null
$ws.Cells.Item(107,5).PasteSpecial(-4122) | Out-Null

# ---- Row 108 ----
$ws.Cells.Item(108,1).Value = 107
$ws.Cells.Item(108,2).Value = 'ecuador'
$ws.Cells.Item(108,3).Value = 'liga-pro'
$ws.Cells.Item(108,4).NumberFormat = "@"
$ws.Cells.Item(108,4).Value = '2023'
$ws.Cells.Item(108,5).Value = 45256
$ws.Cells.Item(108,6).Value = 'Ind. del Valle'
$ws.Cells.Item(108,7).Value = 2
$ws.Cells.Item(108,8).Value = 'Orense'
$ws.Cells.Item(108,9).Value = 2
$ws.Cells.Item(108,10).Value = 1.36
$ws.Cells.Item(108,11).Value = '19/11/2023 00:12'
$ws.Cells.Item(108,12).Value = 1.37
$ws.Cells.Item(108,13).Value = '25/11/2023 23:42'
$ws.Cells.Item(108,14).Value = 4.83
$ws.Cells.Item(108,15).Value = '19/11/2023 00:12'
$ws.Cells.Item(108,16).Value = 4.78
$ws.Cells.Item(108,17).Value = '25/11/2023 23:42'
$ws.Cells.Item(108,18).Value = 8.51
$ws.Cells.Item(108,19).Value = '19/11/2023 00:12'
$ws.Cells.Item(108,20).Value = 8.99
$ws.Cells.Item(108,21).Value = '25/11/2023 23:42'
$ws.Cells.Item(108,22).Value = 'https://www.betexplorer.com/football/ecuador/liga-pro/independiente-del-valle-orense/8p4ERGmk/'
$ws.Cells.Item($lastDataRow,1).Copy() | Out-Null
$ws.Cells.Item(108,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item($lastDataRow,5).Copy() | Out-Null
$ws.Cells.Item(108,5).PasteSpecial(-4122) | Out-Null

# ---- Row 109 ----
$ws.Cells.Item(109,1).Value = 108
$ws.Cells.Item(109,2).Value = 'ecuador'
$ws.Cells.Item(109,3).Value = 'liga-pro'
$ws.Cells.Item(109,4).NumberFormat = "@"
$ws.Cells.Item(109,4).Value = '2023'
$ws.Cells.Item(109,5).Value = 45256
$ws.Cells.Item(109,6).Value = 'Mushuc Runa'
$ws.Cells.Item(109,7).Value = 0
$ws.Cells.Item(109,8).Value = 'U. Catolica'
$ws.Cells.Item(109,9).Value = 2
$ws.Cells.Item(109,10).Value = 3.66
$ws.Cells.Item(109,11).Value = '19/11/2023 00:12'
$ws.Cells.Item(109,12).Value = 3.6
$ws.Cells.Item(109,13).Value = '25/11/2023 23:51'
$ws.Cells.Item(109,14).Value = 3.54
$ws.Cells.Item(109,15).Value = '19/11/2023 00:12'
$ws.Cells.Item(109,16).Value = 3.46
$ws.Cells.Item(109,17).Value = '25/11/2023 23:51'
$ws.Cells.Item(109,18).Value = 1.93
$ws.Cells.Item(109,19).Value = '19/11/2023 00:12'
$ws.Cells.Item(109,20).Value = 2.09
$ws.Cells.Item(109,21).Value = '25/11/2023 23:45'
$ws.Cells.Item(109,22).Value = 'https://www.betexplorer.com/football/ecuador/liga-pro/mushuc-runa-u-catolica/fwOwYIAL/'
$ws.Cells.Item($lastDataRow,1).Copy() | Out-Null
$ws.Cells.Item(109,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item($lastDataRow,5).Copy() | Out-Null
$ws.Cells.Item(109,5).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false
